$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project details")

$ws.Range("A40").Value = "Solar Project DC"
$ws.Range("A41").Value = "Solar Project DC"
$ws.Range("A42").Value = "Solar Project DD"
$ws.Range("A43").Value = "Solar Project AD"
$ws.Range("A44").Value = "Solar Project AD"
$ws.Range("A45").Value = "Solar Project DA"
$ws.Range("A46").Value = "Solar Project DC"
$ws.Range("A47").Value = "Solar Project CB"
$ws.Range("A48").Value = "Solar Project DD"
$ws.Range("A49").Value = "Solar Project BC"
$ws.Range("A50").Value = "Solar Project AD"
$ws.Range("A51").Value = "Solar Project BD"
$ws.Range("A52").Value = "Solar Project AC"
$ws.Range("A53").Value = "Solar Project AD"
$ws.Range("A54").Value = "Solar Project CA"
$ws.Range("A55").Value = "Solar Project BC"
$ws.Range("A56").Value = "Solar Project AB"
$ws.Range("A57").Value = "Solar Project BA"
$ws.Range("A58").Value = "Solar Project AA"
$ws.Range("A59").Value = "Solar Project BC"
$ws.Range("A60").Value = "Solar Project CC"
$ws.Range("A61").Value = "Solar Project BB"
$ws.Range("A62").Value = "Solar Project DA"
$ws.Range("A63").Value = "Solar Project DB"
$ws.Range("A64").Value = "Solar Project BD"
$ws.Range("A65").Value = "Solar Project DB"
$ws.Range("A66").Value = "Solar Project BC"
